$d = $word.ActiveDocument

$firstLineText  = "Mi primera linea"
$secondLineText = "Segunda linea"

# Locate the paragraph that holds the first line of text.
$targetIndex = 1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    $candidateText = $candidate.Range.Text.TrimEnd([char]13, [char]7)
    if ($candidateText -eq $firstLineText) {
        $targetIndex = $i
        break
    }
}
$target = $d.Paragraphs.Item($targetIndex)

# The existing "_GoBack" bookmark currently sits at the end of this paragraph;
# remove it so it can be recreated after the new second line further down.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Collapse to just before the paragraph mark (i.e. right after "Mi primera linea")
# and split the paragraph there, creating a new, second paragraph.
$r = $target.Range
[void]$r.MoveEnd(1, -1)    # wdCharacter: exclude the trailing paragraph mark
[void]$r.Collapse(0)       # wdCollapseEnd
$r.InsertParagraphAfter()

# Type the second line of text into the freshly created paragraph, followed by a
# one-character placeholder. The placeholder lets us anchor a genuinely zero-width
# bookmark right after the real text (anchoring directly at a paragraph's end
# position turns it into a bookmark spanning the wrong content in this runtime),
# after which the placeholder character is deleted again.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newRange = $newPara.Range
[void]$newRange.MoveEnd(1, -1)
$newRange.InsertAfter($secondLineText + "X")

$bmPos = $newPara.Range.Start + $secondLineText.Length
$bmRange = $d.Range($bmPos, $bmPos)
[void]$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($bmPos, $bmPos + 1)
$placeholder.Delete()
